{"js": "// Replace the date and the division-problem text for each table cell.\n// Every \"old\" string below is unique within the document, so a plain\n// search-and-replace on context.document.body is unambiguous.\nconst replacements = [\n  [\"2023-10-08 Sunday\", \"2023-10-09 Monday\"],\n  [\"80\u00f72=40, 0\", \"65\u00f74=16, 1\"],\n  [\"14\u00f72=7, 0\", \"18\u00f73=6, 0\"],\n  [\"61\u00f74=15, 1\", \"99\u00f76=16, 3\"],\n  [\"34\u00f78=4, 2\", \"70\u00f74=17, 2\"],\n  [\"28\u00f76=4, 4\", \"68\u00f74=17, 0\"],\n  [\"59\u00f75=11, 4\", \"48\u00f76=8, 0\"],\n  [\"89\u00f74=22, 1\", \"12\u00f76=2, 0\"],\n  [\"40\u00f79=4, 4\", \"55\u00f76=9, 1\"],\n  [\"92\u00f75=18, 2\", \"99\u00f77=14, 1\"],\n  [\"81\u00f72=40, 1\", \"24\u00f74=6, 0\"],\n  [\"70\u00f73=23, 1\", \"51\u00f79=5, 6\"],\n  [\"82\u00f76=13, 4\", \"99\u00f72=49, 1\"],\n  [\"95\u00f79=10, 5\", \"56\u00f74=14, 0\"],\n  [\"80\u00f79=8, 8\", \"28\u00f74=7, 0\"],\n  [\"17\u00f73=5, 2\", \"41\u00f79=4, 5\"],\n  [\"33\u00f75=6, 3\", \"98\u00f75=19, 3\"],\n  [\"59\u00f78=7, 3\", \"98\u00f78=12, 2\"],\n  [\"76\u00f74=19, 0\", \"98\u00f72=49, 0\"],\n  [\"60\u00f79=6, 6\", \"81\u00f78=10, 1\"],\n  [\"24\u00f77=3, 3\", \"43\u00f76=7, 1\"],\n  [\"13\u00f78=1, 5\", \"91\u00f77=13, 0\"],\n  [\"26\u00f77=3, 5\", \"68\u00f78=8, 4\"],\n  [\"42\u00f77=6, 0\", \"23\u00f77=3, 2\"],\n  [\"87\u00f79=9, 6\", \"76\u00f76=12, 4\"],\n  [\"84\u00f78=10, 4\", \"39\u00f74=9, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the division-problem text for each table cell.\n# Every \"old\" string below is unique within the document, so a plain\n# Find/Replace against the whole document body is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2023-10-08 Sunday\", \"2023-10-09 Monday\"),\n  @(\"80\u00f72=40, 0\", \"65\u00f74=16, 1\"),\n  @(\"14\u00f72=7, 0\", \"18\u00f73=6, 0\"),\n  @(\"61\u00f74=15, 1\", \"99\u00f76=16, 3\"),\n  @(\"34\u00f78=4, 2\", \"70\u00f74=17, 2\"),\n  @(\"28\u00f76=4, 4\", \"68\u00f74=17, 0\"),\n  @(\"59\u00f75=11, 4\", \"48\u00f76=8, 0\"),\n  @(\"89\u00f74=22, 1\", \"12\u00f76=2, 0\"),\n  @(\"40\u00f79=4, 4\", \"55\u00f76=9, 1\"),\n  @(\"92\u00f75=18, 2\", \"99\u00f77=14, 1\"),\n  @(\"81\u00f72=40, 1\", \"24\u00f74=6, 0\"),\n  @(\"70\u00f73=23, 1\", \"51\u00f79=5, 6\"),\n  @(\"82\u00f76=13, 4\", \"99\u00f72=49, 1\"),\n  @(\"95\u00f79=10, 5\", \"56\u00f74=14, 0\"),\n  @(\"80\u00f79=8, 8\", \"28\u00f74=7, 0\"),\n  @(\"17\u00f73=5, 2\", \"41\u00f79=4, 5\"),\n  @(\"33\u00f75=6, 3\", \"98\u00f75=19, 3\"),\n  @(\"59\u00f78=7, 3\", \"98\u00f78=12, 2\"),\n  @(\"76\u00f74=19, 0\", \"98\u00f72=49, 0\"),\n  @(\"60\u00f79=6, 6\", \"81\u00f78=10, 1\"),\n  @(\"24\u00f77=3, 3\", \"43\u00f76=7, 1\"),\n  @(\"13\u00f78=1, 5\", \"91\u00f77=13, 0\"),\n  @(\"26\u00f77=3, 5\", \"68\u00f78=8, 4\"),\n  @(\"42\u00f77=6, 0\", \"23\u00f77=3, 2\"),\n  @(\"87\u00f79=9, 6\", \"76\u00f76=12, 4\"),\n  @(\"84\u00f78=10, 4\", \"39\u00f74=9, 3\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
